$wb = $excel.ActiveWorkbook

# Row 32 on sheet ALC (diff @ -2197)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 435.55554
$ws.Range("I32").Value = 714.2857
$ws.Range("J32").Value = 258.18182
$ws.Range("K32").Value = 714.2857
$ws.Range("L32").Value = 258.18182
$ws.Range("M32").Value = -388.2857
$ws.Range("N32").Value = -910.18182

# Row 121 on sheet ALC (diff @ -6690)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H121").Value = 937.1429000000001
$ws.Range("J121").Value = 1001.1539
$ws.Range("L121").Value = 3003.4617
$ws.Range("N121").Value = -6497.4617

# Row 137 on sheet ALC (diff @ -7495)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 29038.277
$ws.Range("I137").Value = 877.5238000000001
$ws.Range("J137").Value = 68463.336
$ws.Range("K137").Value = 2632.5714
$ws.Range("L137").Value = 205390.008
$ws.Range("M137").Value = -82.57140000000027
$ws.Range("N137").Value = -210490.008

# Row 32 on sheet ARM (diff @ -9325)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5653.48
$ws.Range("I32").Value = 4391.1904
$ws.Range("J32").Value = 12280.5
$ws.Range("K32").Value = 4391.1904
$ws.Range("L32").Value = 12280.5
$ws.Range("M32").Value = -4104.1904
$ws.Range("N32").Value = -12854.5

# Row 61 on sheet ARM (diff @ -10749)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1865.9166
$ws.Range("I61").Value = 1647.6578
$ws.Range("J61").Value = 2695.3
$ws.Range("K61").Value = 1647.6578
$ws.Range("L61").Value = 2695.3
$ws.Range("M61").Value = -1435.6578
$ws.Range("N61").Value = -3119.3

# Row 74 on sheet ARM (diff @ -11383)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 229514.4
$ws.Range("I74").Value = 2560.9285
$ws.Range("J74").Value = 626683
$ws.Range("K74").Value = 2560.9285
$ws.Range("L74").Value = 626683
$ws.Range("M74").Value = -1686.9285
$ws.Range("N74").Value = -628431

# Row 77 on sheet ARM (diff @ -11527)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 229514.4
$ws.Range("I77").Value = 2560.9285
$ws.Range("J77").Value = 626683
$ws.Range("K77").Value = 12804.6425
$ws.Range("L77").Value = 3133415
$ws.Range("M77").Value = -8436.6425
$ws.Range("N77").Value = -3142151

# Row 120 on sheet ARM (diff @ -13631)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H120").Value = 30242.666
$ws.Range("J120").Value = 30242.666
$ws.Range("L120").Value = 30242.666
$ws.Range("N120").Value = -39918.666

# Row 132 on sheet ARM (diff @ -14222)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 18448
$ws.Range("I132").Value = 27771.578
$ws.Range("K132").Value = 83314.734
$ws.Range("M132").Value = -80784.734

# Row 136 on sheet ARM (diff @ -14421)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1865.9166
$ws.Range("I136").Value = 1647.6578
$ws.Range("J136").Value = 2695.3
$ws.Range("K136").Value = 4942.9734
$ws.Range("L136").Value = 8085.900000000001
$ws.Range("M136").Value = -2392.9734
$ws.Range("N136").Value = -13185.9

# Row 86 on sheet BSM (diff @ -18940)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 5558082.5
$ws.Range("I86").Value = 8697935
$ws.Range("J86").Value = 2959.3845
$ws.Range("K86").Value = 8697935
$ws.Range("L86").Value = 2959.3845
$ws.Range("M86").Value = -8696812
$ws.Range("N86").Value = -5205.3845

# Row 89 on sheet BSM (diff @ -19090)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 5558082.5
$ws.Range("I89").Value = 8697935
$ws.Range("J89").Value = 2959.3845
$ws.Range("K89").Value = 43489675
$ws.Range("L89").Value = 14796.9225
$ws.Range("M89").Value = -43484059
$ws.Range("N89").Value = -26028.9225

# Row 134 on sheet BSM (diff @ -21304)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 346171.22
$ws.Range("I134").Value = 417532.53
$ws.Range("J134").Value = 3636.8
$ws.Range("K134").Value = 1252597.59
$ws.Range("L134").Value = 10910.4
$ws.Range("M134").Value = -1250062.59
$ws.Range("N134").Value = -15980.4

# Row 31 on sheet CRP (diff @ -23229)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1835.322
$ws.Range("I31").Value = 1194.8286
$ws.Range("J31").Value = 2769.375
$ws.Range("K31").Value = 1194.8286
$ws.Range("L31").Value = 2769.375
$ws.Range("M31").Value = -899.8286000000001
$ws.Range("N31").Value = -3359.375

# Row 34 on sheet CRP (diff @ -23382)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1835.322
$ws.Range("I34").Value = 1194.8286
$ws.Range("J34").Value = 2769.375
$ws.Range("K34").Value = 1194.8286
$ws.Range("L34").Value = 2769.375
$ws.Range("M34").Value = -992.8286000000001
$ws.Range("N34").Value = -3173.375

# Row 99 on sheet CRP (diff @ -26555)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 5515.2
$ws.Range("I99").Value = 5736
$ws.Range("J99").Value = 5000
$ws.Range("K99").Value = 5736
$ws.Range("L99").Value = 5000
$ws.Range("M99").Value = -4238
$ws.Range("N99").Value = -7996

# Row 126 on sheet CRP (diff @ -27866)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 5515.2
$ws.Range("I126").Value = 5736
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 17208
$ws.Range("L126").Value = 15000
$ws.Range("M126").Value = -14738
$ws.Range("N126").Value = -19940

# Row 5 on sheet CUL (diff @ -28903)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 782.95746
$ws.Range("I5").Value = 457.96
$ws.Range("K5").Value = 1373.88
$ws.Range("M5").Value = -1261.88

# Row 26 on sheet CUL (diff @ -29971)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 533.3333
$ws.Range("I26").Value = 475
$ws.Range("J26").Value = 650
$ws.Range("K26").Value = 1425
$ws.Range("L26").Value = 1950
$ws.Range("M26").Value = -1137
$ws.Range("N26").Value = -2526

# Row 50 on sheet CUL (diff @ -31186)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 164.55556
$ws.Range("I50").Value = 47
$ws.Range("J50").Value = 399.66666
$ws.Range("K50").Value = 141
$ws.Range("L50").Value = 1198.99998
$ws.Range("M50").Value = 340
$ws.Range("N50").Value = -2160.99998

# Row 53 on sheet CUL (diff @ -31342)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H53").Value = 164.55556
$ws.Range("I53").Value = 47
$ws.Range("J53").Value = 399.66666
$ws.Range("K53").Value = 141
$ws.Range("L53").Value = 1198.99998
$ws.Range("M53").Value = 340
$ws.Range("N53").Value = -2160.99998

# Row 122 on sheet CUL (diff @ -34864)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 844.0968
$ws.Range("I122").Value = 457
$ws.Range("J122").Value = 1314.1428
$ws.Range("K122").Value = 4113
$ws.Range("L122").Value = 11827.2852
$ws.Range("M122").Value = -1663
$ws.Range("N122").Value = -16727.2852

# Row 131 on sheet CUL (diff @ -35323)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1961474.9
$ws.Range("J131").Value = 981.1111
$ws.Range("L131").Value = 2943.3333
$ws.Range("N131").Value = -13023.3333

# Row 132 on sheet CUL (diff @ -35375)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1144.625
$ws.Range("I132").Value = 784.6923
$ws.Range("J132").Value = 1570
$ws.Range("K132").Value = 7062.2307
$ws.Range("L132").Value = 14130
$ws.Range("M132").Value = -4532.2307
$ws.Range("N132").Value = -19190

# Row 135 on sheet CUL (diff @ -35528)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 782.95746
$ws.Range("I135").Value = 457.96
$ws.Range("K135").Value = 4121.639999999999
$ws.Range("M135").Value = -1586.639999999999

# Row 41 on sheet GSM (diff @ -37909)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 3051
$ws.Range("I41").Value = 3051
$ws.Range("K41").Value = 3051
$ws.Range("M41").Value = -2696

# Row 132 on sheet GSM (diff @ -42296)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2644.919
$ws.Range("I132").Value = 1898.5217
$ws.Range("J132").Value = 3871.1428
$ws.Range("K132").Value = 5695.5651
$ws.Range("L132").Value = 11613.4284
$ws.Range("M132").Value = -3165.5651
$ws.Range("N132").Value = -16673.4284

# Row 110 on sheet LTW (diff @ -48142)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H110").Value = 21411
$ws.Range("J110").Value = 21411
$ws.Range("L110").Value = 21411
$ws.Range("N110").Value = -29591

# Row 136 on sheet LTW (diff @ -49407)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 19457.125
$ws.Range("J136").Value = 6312.25
$ws.Range("L136").Value = 18936.75
$ws.Range("N136").Value = -24036.75
